$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.083773612976074
$ws.Range("B1").Value = 1.809671759605408
$ws.Range("C1").Value = 5.327514171600342
$ws.Range("D1").Value = 0.7834926843643188
$ws.Range("E1").Value = 0.555210292339325
